{"js": "// The document contains literal escape-style placeholders of the form\n// \"<u+XXXX>\" / \"<U+XXXX>\" (produced by an upstream HTML/unicode-escape\n// export step) sitting inside otherwise normal run text, e.g.\n//   \"...mellitus<u+2013>a systematic review\"\n// This script decodes every such placeholder back into the real Unicode\n// character it represents, e.g. \"...mellitus\\u2013a systematic review\".\n//\n// Special case: \"<u+00A0>\" (NBSP) was rendered back as a plain space in\n// the target document, not an actual U+00A0 character, so it is handled\n// separately.\n\nconst body = context.document.body;\n\n// code point (uppercase hex) -> replacement string\nconst decodeMap = {\n  \"2013\": \"\\u2013\", // EN DASH\n  \"2014\": \"\\u2014\", // EM DASH\n  \"00D8\": \"\\u00D8\", // \u00d8\n  \"00F8\": \"\\u00F8\", // \u00f8\n  \"00E6\": \"\\u00E6\", // \u00e6\n  \"00E1\": \"\\u00E1\", // \u00e1\n  \"00A0\": \" \"       // NBSP placeholder -> plain space (per target content)\n};\n\nfor (const code of Object.keys(decodeMap)) {\n  const token = \"<u+\" + code + \">\";\n  const replacement = decodeMap[code];\n\n  const results = body.search(token, { matchCase: false, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains literal escape-style placeholders of the form\n# \"<u+XXXX>\" / \"<U+XXXX>\" (produced by an upstream HTML/unicode-escape\n# export step) sitting inside otherwise normal run text, e.g.\n#   \"...mellitus<u+2013>a systematic review\"\n# This script decodes every such placeholder back into the real Unicode\n# character it represents, e.g. \"...mellitus\" + (EN DASH) + \"a systematic review\".\n#\n# Special case: \"<u+00A0>\" (NBSP) was rendered back as a plain space in\n# the target document, not an actual U+00A0 character, so it is handled\n# separately.\n\n$d = $word.ActiveDocument\n\n$decodeMap = @(\n  @{code = \"2013\"; repl = [char]0x2013},  # EN DASH\n  @{code = \"2014\"; repl = [char]0x2014},  # EM DASH\n  @{code = \"00D8\"; repl = [char]0x00D8},  # \u00d8\n  @{code = \"00F8\"; repl = [char]0x00F8},  # \u00f8\n  @{code = \"00E6\"; repl = [char]0x00E6},  # \u00e6\n  @{code = \"00E1\"; repl = [char]0x00E1},  # \u00e1\n  @{code = \"00A0\"; repl = \" \"}            # NBSP placeholder -> plain space\n)\n\nforeach ($entry in $decodeMap) {\n  $token = \"<u+\" + $entry.code + \">\"\n\n  $find = $d.Content.Find\n  $find.Text = $token\n  $find.MatchCase = $false\n  $find.MatchWildcards = $false\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Replacement.Text = $entry.repl\n  $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
